# Applies the diff: updates the PSRR-related rows of the specs table.
#   Row "Min PSRR" / "> 5 dB"              -> "PSRR at 1 MHz" / "> 30 dB"
#   Row "PSRR (for freq < 10M)" / "> 30 dB" -> "PSRR at DC"     / "> 40 dB"

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-CellRuns($cell, [string]$innerRunsXml) {
    $rng = $cell.Range
    # Exclude the trailing end-of-cell marker so only the paragraph content is replaced.
    $target = $d.Range($rng.Start, $rng.End - 1)
    $xml = $pkgHeader + '<w:p>' + $innerRunsXml + '</w:p>' + $pkgFooter
    $target.InsertXML($xml)
}

# Row: "Min PSRR" -> "PSRR" + " at 1 MHz"
Set-CellRuns $t.Cell(9, 1) '<w:r><w:t>PSRR</w:t></w:r><w:r><w:t xml:space="preserve"> at 1 MHz</w:t></w:r>'

# Row: "> 5 dB" -> "> " + "30" + " dB"
Set-CellRuns $t.Cell(9, 2) '<w:r><w:t xml:space="preserve">&gt; </w:t></w:r><w:r><w:t>30</w:t></w:r><w:r><w:t xml:space="preserve"> dB</w:t></w:r>'

# Row: "PSRR (for freq < 10M)" -> "PSRR " + "at DC"
Set-CellRuns $t.Cell(10, 1) '<w:r><w:t xml:space="preserve">PSRR </w:t></w:r><w:r><w:t>at DC</w:t></w:r>'

# Row: "> 30 dB" -> "> " + "4" + "0 dB"
Set-CellRuns $t.Cell(10, 2) '<w:r><w:t xml:space="preserve">&gt; </w:t></w:r><w:r><w:t>4</w:t></w:r><w:r><w:t>0 dB</w:t></w:r>'

Write-Host "Row9 C1: [" $t.Cell(9,1).Range.Text "]"
Write-Host "Row9 C2: [" $t.Cell(9,2).Range.Text "]"
Write-Host "Row10 C1: [" $t.Cell(10,1).Range.Text "]"
Write-Host "Row10 C2: [" $t.Cell(10,2).Range.Text "]"
